$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Novos resultados (Foto 8: Fundo / Folha / Doença)
$ws.Range("H25").Value = 613325
$ws.Range("H26").Value = 298611
$ws.Range("H27").Formula = "=921557 - SUM(H25:H26)"

# Atualiza a celula ativa/selecao, conforme o arquivo final
$ws.Range("I24").Select()
